$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.015.54'
$ws.Range('E2').Value = '  +1.21%  '
$ws.Range('D3').Value = '2.341.02'
$ws.Range('E3').Value = '  +0.96%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '518.53'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.30%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '135.85'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.07%  '
$ws.Range('E7').Value = '  +0.19%  '
$ws.Range('D9').Value = '2.350.93'
$ws.Range('E9').Value = '  +0.42%  '
$ws.Range('E10').Value = '  -0.54%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.41'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +5.04%  '
$ws.Range('E12').Value = '  -1.39%  '
$ws.Range('E13').Value = '  -0.05%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '23.96'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.50%  '
$ws.Range('D15').Value = '2.757.68'
$ws.Range('E15').Value = '  +0.18%  '
$ws.Range('D16').Value = '56.988.08'
$ws.Range('E16').Value = '  +1.07%  '
$ws.Range('E17').Value = '  -0.34%  '
$ws.Range('D18').Value = '2.348.98'
$ws.Range('E18').Value = '  +0.47%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.61'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.63%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '327.12'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.49%  '
$ws.Range('E21').Value = '  -0.35%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.77'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.17%  '
$ws.Range('E23').Value = '  -0.28%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '61.26'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.52%  '
$ws.Range('E25').Value = '  +5.06%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.996'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.34%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.04'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +5.28%  '
$ws.Range('E28').Value = '  +9.42%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '170.15'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.11%  '
$ws.Range('D30').Value = '0.0₃0742'
$ws.Range('E30').Value = '  +2.68%  '
$ws.Range('E31').Value = '  +1.19%  '
$ws.Range('E32').Value = '  +0.15%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '18.56'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.13%  '
$ws.Range('E34').Value = '  +0.02%  '
$ws.Range('E35').Value = '  +0.22%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.28'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.06%  '
$ws.Range('E37').Value = '  -1.67%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.02'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.76%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '38.46'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.96%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '149.73'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +7.86%  '
$ws.Range('E42').Value = '  +0.17%  '
$ws.Range('E43').Value = '  +0.68%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '280.32'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +4.58%  '
$ws.Range('E45').Value = '  +1.54%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0935'
$ws.Range('D46').Style = 'Normal'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0507'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.64%  '
$ws.Range('E48').Value = '  +1.47%  '
$ws.Range('E49').Value = '  +2.05%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '18.02'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +6.03%  '
$ws.Range('E51').Value = '  -0.09%  '
